# Apply the MapNpcData "list column" addition edit.
#
# Summary of the target change:
#  - Sheet "MapNpcInfoData": the header text for columns C2/E2 is pluralized
#    ("ShowRequirementValue" -> "ShowRequirementValues",
#     "HideRequirementValue" -> "HideRequirementValues").
#  - Sheet "MapNpcMenuData": two new "Values" columns (list<int>) are
#    inserted - one right after "ShowRequirement" (old column B) and one
#    right after "HideRequirement" (old column C, which becomes column D
#    after the first insert). The rest of the columns shift right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) MapNpcInfoData: rename the two header cells (simple text edit).
#    Update E2 before C2 so the shared-string table keeps the same
#    slot assignment the original authors' edit produced.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MapNpcInfoData")
$ws1.Range("E2").Value = "HideRequirementValues"
$ws1.Range("C2").Value = "ShowRequirementValues"
[void]$ws1.Range("E6").Select()

# ---------------------------------------------------------------------
# 2) MapNpcMenuData: insert the two new columns.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MapNpcMenuData")

# Remember the column widths we need to propagate to the new columns
# before the column letters shift around.
$widthShow = $ws2.Range("B1").ColumnWidth()   # ShowRequirement column width
$widthHide = $ws2.Range("C1").ColumnWidth()   # HideRequirement column width

# Insert a new column right after "ShowRequirement" (old B) -- it
# inherits ShowRequirement's column formatting automatically.
$ws2.Columns("C:C").Insert()
$ws2.Range("C1").ColumnWidth = $widthShow

# Insert a new column right after "HideRequirement" (old C, now shifted
# to D). It first inherits HideRequirement's column formatting...
$ws2.Columns("E:E").Insert()
$ws2.Range("E1").ColumnWidth = $widthHide

# ...but the header rows (1-2) of the new "Values" column should look
# like the "ShowRequirementValues" column (C), not like "HideRequirement"
# (D), so restore that formatting now.
[void]$ws2.Range("C1:C2").Copy()
[void]$ws2.Range("E1:E2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Finally fill in the header text for the two new columns. Set E's
# (Hide) text before C's (Show) text, matching the shared-string slot
# order used above.
$ws2.Range("E1").Value = "list<int>"
$ws2.Range("E2").Value = "HideRequirementValues"
$ws2.Range("C1").Value = "list<int>"
$ws2.Range("C2").Value = "ShowRequirementValues"

[void]$ws2.Range("C20").Select()

# Restore MapNpcInfoData as the active sheet/tab.
[void]$ws1.Activate()
